$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

# Only columns A (statut) and B (statut_label) ever hold the
# "⬛"/"🟩"/"noir" text values being renamed here, so restrict the scan
# to those two columns (avoids false matches against boolean cells,
# where PowerShell's "-eq" coerces the string operand to the left
# operand's [bool] type).
for ($r = 2; $r -le $rows; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value()
    if ($valA -is [string]) {
        if ($valA -eq "⬛") {
            $cellA.Value = "📘"
        } elseif ($valA -eq "🟩") {
            $cellA.Value = "📗"
        }
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value()
    if ($valB -is [string]) {
        if ($valB -eq "noir") {
            $cellB.Value = "bleu"
        }
    }
}
